$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "44.057.44"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "2.360.41"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.93"
$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.40"
$ws.Range("E7").Value = "  +2.57%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  +11.53%  "

$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.26"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.25"
$ws.Range("E12").Value = "  +10.80%  "

$ws.Range("E13").Value = "  +9.88%  "

$ws.Range("D15").Value = "2.712.23"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.64"
$ws.Range("E16").Value = "  -0.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.900"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").Value = "2.354.25"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").Value = "43.959.86"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("E21").Value = "  +5.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.04"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.63"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +25.46%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.50"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.66"
$ws.Range("E27").Value = "  -2.40%  "

$ws.Range("E28").Value = "  +3.01%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.78"
$ws.Range("E30").Value = "  +1.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.80"
$ws.Range("E31").Value = "  +1.92%  "

$ws.Range("E32").Value = "  -1.95%  "

$ws.Range("E33").Value = "  +3.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0773"
$ws.Range("E34").Value = "  +7.68%  "

$ws.Range("E35").Value = "  +1.74%  "

$ws.Range("E36").Value = "  +4.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").Value = "  -3.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0281"
$ws.Range("E40").Value = "  +5.40%  "

$ws.Range("E41").Value = "  +15.32%  "

$ws.Range("E42").Value = "  +14.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.12"
$ws.Range("E43").Value = "  +3.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.13"
$ws.Range("E44").Value = "  -1.41%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.76"
$ws.Range("E46").Value = "  +6.84%  "

$ws.Range("E47").Value = "  +9.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "58.11"
$ws.Range("E48").Value = "  +10.98%  "

$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("E50").Value = "  +0.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.49"
$ws.Range("E51").Value = "  +2.20%  "
